$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in measured Vout (dBvrms) values in column E (rows 3-18) ---
$ws.Range("E3").Value = -128.05000000000001
$ws.Range("E4").Value = -111.813
$ws.Range("E5").Value = -56.6
$ws.Range("E6").Value = -49
$ws.Range("E7").Value = -41.5
$ws.Range("E8").Value = -34.5
$ws.Range("E9").Value = -27.7
$ws.Range("E10").Value = -21.1
$ws.Range("E11").Value = -14.8
$ws.Range("E12").Value = -9.6999999999999993
$ws.Range("E13").Value = -7.4
$ws.Range("E14").Value = -6.5
$ws.Range("E15").Value = -6.3
$ws.Range("E16").Value = -6.3
$ws.Range("E17").Value = -6.2
$ws.Range("E18").Value = -6.2

# --- New report block starting row 23 ---
$ws.Range("D23").Value = "Vout(V)"
$ws.Range("E23").Value = "vout*vin"
$ws.Range("F23").Value = "DB"

$ws.Range("D24").Formula = "=10^(E3/20)"
$ws.Range("D25").Formula = "=10^(E4/20)"
$ws.Range("D26").Formula = "=10^(E5/20)"
$ws.Range("D27").Formula = "=10^(E6/20)"
$ws.Range("D28").Formula = "=10^(E7/20)"
$ws.Range("D29").Formula = "=10^(E8/20)"
$ws.Range("D30").Formula = "=10^(E9/20)"
$ws.Range("D31").Formula = "=10^(E10/20)"
$ws.Range("D32").Formula = "=10^(E11/20)"
$ws.Range("D33").Formula = "=10^(E12/20)"
$ws.Range("D34").Formula = "=10^(E13/20)"
$ws.Range("D35").Formula = "=10^(E14/20)"
$ws.Range("D36").Formula = "=10^(E15/20)"
$ws.Range("D37").Formula = "=10^(E16/20)"
$ws.Range("D38").Formula = "=10^(E17/20)"
$ws.Range("D39").Formula = "=10^(E18/20)"

$ws.Range("E24").Formula = "=D24*C4"
$ws.Range("E25").Formula = "=D25*C5"
$ws.Range("E26").Formula = "=D26*C6"
$ws.Range("E27").Formula = "=D27*C7"
$ws.Range("E28").Formula = "=D28*C8"
$ws.Range("E29").Formula = "=D29*C9"
$ws.Range("E30").Formula = "=D30*C10"
$ws.Range("E31").Formula = "=D31*C11"
$ws.Range("E32").Formula = "=D32*C12"
$ws.Range("E33").Formula = "=D33*C13"
$ws.Range("E34").Formula = "=D34*C14"
$ws.Range("E35").Formula = "=D35*C15"
$ws.Range("E36").Formula = "=D36*C16"
$ws.Range("E37").Formula = "=D37*C17"
$ws.Range("E38").Formula = "=D38*C18"
$ws.Range("E39").Formula = "=D39*C19"

$ws.Range("F24").Formula = "= 20*LOG10(E24)"
$ws.Range("F25").Formula = "= 20*LOG10(E25)"
$ws.Range("F26").Formula = "= 20*LOG10(E26)"
$ws.Range("F27").Formula = "= 20*LOG10(E27)"
$ws.Range("F28").Formula = "= 20*LOG10(E28)"
$ws.Range("F29").Formula = "= 20*LOG10(E29)"
$ws.Range("F30").Formula = "= 20*LOG10(E30)"
$ws.Range("F31").Formula = "= 20*LOG10(E31)"
$ws.Range("F32").Formula = "= 20*LOG10(E32)"
$ws.Range("F33").Formula = "= 20*LOG10(E33)"
$ws.Range("F34").Formula = "= 20*LOG10(E34)"
$ws.Range("F35").Formula = "= 20*LOG10(E35)"
$ws.Range("F36").Formula = "= 20*LOG10(E36)"
$ws.Range("F37").Formula = "= 20*LOG10(E37)"
$ws.Range("F38").Formula = "= 20*LOG10(E38)"
$ws.Range("F39").Formula = "= 20*LOG10(E39)"

# --- Resize chart 2 to cover the extended data range ---
$chart2 = $ws.ChartObjects(2)
$chart2.Width = 1001.7201771653542
$chart2.Height = 877.7283464566929

# --- View state tweaks ---
$excel.ActiveWindow.Zoom = 175
$ws.Range("F38").Select()
